# Closed the Europa budget and changed the notation of the coding in the
# original as well to keep it compatible.
#
# Row 21 ("Modulation/coding type") previously held the modulation type
# "FSK" in column B with a numeric "8" (presumably an order, bits, or
# similar) repeated across C:H for every spacecraft budget. The edit
# replaces the modulation description with the single notation "8FSK"
# used consistently for every column (C:H), and sets column B (which used
# to hold the modulation name) to the generic "-" placeholder used
# elsewhere in the sheet, since the notation now lives directly in the
# data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "-"
$ws.Range("C21:H21").Value = "8FSK"

# Update the selection to reflect where the author was last working.
$ws.Range("I28").Select() | Out-Null
